$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: duplicate applicant (Kelly Marinduque) moving to a new position (System Administrator), referred by Hart
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A3").Value = 42915
$ws.Range("B3").Value = "Kelly Marinduque"
$ws.Range("C3").Value = "Hart"
$ws.Range("D3").Value = "System Administrator"
$ws.Range("E3").Value = 565434

# Row 4: new applicant (Leah Mahusay) referred, applying for Fiber Technician
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A4").Value = 42915
$ws.Range("B4").Value = "Leah Mahusay"
$ws.Range("C4").Value = "Referral"
$ws.Range("D4").Value = "Fiber Technician"
$ws.Range("E4").Value = 124325434

$ws.Range("E3").Select()
